$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text representation (values such as "337.15"
# would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.964.00'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '1.760.72'
$ws.Range('E3').Value = '  -3.44%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '337.15'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('D6').Value = '0.9991'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.3779'
$ws.Range('E7').Value = '  -4.51%  '
$ws.Range('D8').Value = '0.3346'
$ws.Range('E8').Value = '  -4.70%  '
$ws.Range('D9').Value = '45.72'
$ws.Range('E9').Value = '  -5.14%  '
$ws.Range('D10').Value = '1.121'
$ws.Range('E10').Value = '  -7.03%  '
$ws.Range('D11').Value = '0.07188'
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').Value = '22.32'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '6.182'
$ws.Range('E14').Value = '  -5.88%  '
$ws.Range('D15').Value = '7.179'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '1.758.00'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('E17').Value = '  -5.29%  '
$ws.Range('D18').Value = '0.06571'
$ws.Range('E18').Value = '  -2.44%  '
$ws.Range('D19').Value = '80.38'
$ws.Range('E19').Value = '  -6.24%  '
$ws.Range('D20').Value = '0.9999'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  -5.45%  '
$ws.Range('D22').Value = '6.278'
$ws.Range('E22').Value = '  -5.14%  '
$ws.Range('D23').Value = '27.978.93'
$ws.Range('E23').Value = '  -0.75%  '
$ws.Range('D24').Value = '11.68'
$ws.Range('E24').Value = '  -8.88%  '
$ws.Range('D25').Value = '2.363'
$ws.Range('E25').Value = '  -1.79%  '
$ws.Range('D26').Value = '152.11'
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('D27').Value = '19.86'
$ws.Range('E27').Value = '  -7.82%  '
$ws.Range('D28').Value = '2.331'
$ws.Range('E28').Value = '  -9.80%  '
$ws.Range('D29').Value = '1.958.59'
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('D30').Value = '1.267'
$ws.Range('E30').Value = '  -15.91%  '
$ws.Range('D31').Value = '131.90'
$ws.Range('E31').Value = '  -3.59%  '
$ws.Range('D32').Value = '4.018'
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('D33').Value = '5.793'
$ws.Range('E33').Value = '  -7.01%  '
$ws.Range('D34').Value = '0.08760'
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').Value = '12.24'
$ws.Range('E35').Value = '  -7.61%  '
$ws.Range('D36').Value = '0.02334'
$ws.Range('E36').Value = '  -4.89%  '
$ws.Range('D37').Value = '0.6579'
$ws.Range('E37').Value = '  -6.25%  '
$ws.Range('D38').Value = '0.06186'
$ws.Range('E38').Value = '  -6.24%  '
$ws.Range('D39').Value = '5.144'
$ws.Range('E39').Value = '  -7.44%  '
$ws.Range('D40').Value = '0.2103'
$ws.Range('E40').Value = '  -6.19%  '
$ws.Range('D41').Value = '1.209'
$ws.Range('E41').Value = '  -5.03%  '
$ws.Range('D42').Value = '1.446'
$ws.Range('E42').Value = '  -10.50%  '
$ws.Range('D43').Value = '7.993'
$ws.Range('E43').Value = '  -6.73%  '
$ws.Range('D44').Value = '0.9993'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '13.71'
$ws.Range('E45').Value = '  -6.57%  '
$ws.Range('D46').Value = '3.823'
$ws.Range('E46').Value = '  -1.90%  '
$ws.Range('D47').Value = '0.6038'
$ws.Range('E47').Value = '  -7.76%  '
$ws.Range('D48').Value = '130.18'
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('D49').Value = '2.009'
$ws.Range('E49').Value = '  -8.16%  '
$ws.Range('D50').Value = '1.184'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('D51').Value = '0.07155'
$ws.Range('E51').Value = '  -1.10%  '
